# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for the data rows (16-23): the "DAIRY LUZ CAUSADO CERVANTES"
# worker record moves to the top of the list (row 16), and the remaining
# "YESENIA ORTEGA RODRIGUEZ" period rows are re-sorted in ascending order
# (2310 .. 2404). The "Valor Mora" amounts (column F) travel together with
# each (document, period) row.
$rows = @(
    @{ B = "CC"; C = "1143355624"; D = "DAIRY LUZ CAUSADO CERVANTES"; E = "1802"; F = 48000 },
    @{ B = "CC"; C = "45694304";   D = "YESENIA ORTEGA RODRIGUEZ";    E = "2310"; F = 48000 },
    @{ B = "CC"; C = "45694304";   D = "YESENIA ORTEGA RODRIGUEZ";    E = "2311"; F = 48000 },
    @{ B = "CC"; C = "45694304";   D = "YESENIA ORTEGA RODRIGUEZ";    E = "2312"; F = 48000 },
    @{ B = "CC"; C = "45694304";   D = "YESENIA ORTEGA RODRIGUEZ";    E = "2401"; F = 48000 },
    @{ B = "CC"; C = "45694304";   D = "YESENIA ORTEGA RODRIGUEZ";    E = "2402"; F = 48000 },
    @{ B = "CC"; C = "45694304";   D = "YESENIA ORTEGA RODRIGUEZ";    E = "2403"; F = 48000 },
    @{ B = "CC"; C = "45694304";   D = "YESENIA ORTEGA RODRIGUEZ";    E = "2404"; F = 36800 }
)

$startRow = 16
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
}

$wb.Save()
